$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
  81.501057082452434,
  78.752642706131084,
  79.492600422832979,
  77.167019027484145,
  83.826638477801268,
  87.315010570824526,
  83.192389006342495,
  83.086680761099359,
  83.192389006342495,
  83.298097251585617,
  86.892177589852011,
  86.786469344608875,
  83.509513742071888,
  83.61522198731501,
  83.509513742071888,
  83.720930232558146,
  87.20930232558139,
  87.20930232558139,
  79.598308668076115,
  72.41014799154334,
  79.704016913319236,
  79.175475687103585,
  80.021141649048616,
  79.281183932346721,
  83.61522198731501,
  83.61522198731501,
  84.143763213530647,
  84.883720930232556,
  84.883720930232556,
  85.306553911205071,
  74.524312896405917,
  80.655391120507403,
  80.866807610993661,
  76.849894291754765,
  76.955602536997887,
  80.655391120507403,
  76.849894291754765,
  76.744186046511629,
  79.704016913319236,
  83.720930232558146,
  82.980972515856237,
  83.509513742071888,
  82.663847780126858,
  82.135306553911207,
  82.135306553911207,
  82.875264270613101,
  77.061310782241009,
  82.980972515856237,
  83.192389006342495
)

for ($i = 0; $i -lt $values.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $values[$i]
}
